# Updated cryptos list on Sun Aug 20 16:35:29 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (D) and "Volume(1h)" (E) columns of the crypto
# table with the latest scraped quotes. Both columns are stored as plain
# text in the sheet (Price sometimes uses more than one "." as a
# thousands separator, e.g. "26.370.54", and Volume is a padded
# "  +0.44%  " style string) so values are written back as text.
#
# Most Price values parse as a genuine number when assigned through
# Range.Value (just like typing them into Excel would), which would
# silently flip the cell from a text cell to a numeric one. To keep
# those cells textual - matching the source data - NumberFormat is set
# to "@" (Text) before the value is written for exactly those cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> new Price (column D) text. Only rows whose Price changed.
$priceUpdates = @{
    2  = "26.370.54"
    3  = "1.687.86"
    5  = "218.42"
    6  = "0.5461"
    8  = "0.2725"
    9  = "0.06463"
    10 = "22.01"
    11 = "0.07686"
    12 = "1.709.97"
    13 = "4.536"
    14 = "0.5809"
    15 = "0.000008380"
    16 = "65.09"
    17 = "26.426.75"
    21 = "190.72"
    22 = "6.245"
    23 = "1.011"
    24 = "149.75"
    25 = "0.1313"
    26 = "7.874"
    28 = "0.06326"
    29 = "1.410"
    30 = "1.326"
    31 = "3.588"
    32 = "3.579"
    33 = "1.685"
    35 = "0.6175"
    37 = "2.719"
    38 = "6.270"
    39 = "1.112.35"
    40 = "0.01626"
    41 = "0.8799"
    43 = "101.42"
    45 = "0.00000000110"
    46 = "57.34"
    47 = "1.013"
    48 = "8.171"
    50 = "0.4305"
    51 = "6.039"
}

# row -> new Volume(1h) (column E) text. Only rows whose Volume changed.
$volumeUpdates = @{
    2  = "  +0.44%  "
    3  = "  +0.02%  "
    5  = "  -0.24%  "
    6  = "  +4.04%  "
    7  = "  +0.61%  "
    8  = "  +1.28%  "
    9  = "  +0.30%  "
    10 = "  -0.44%  "
    11 = "  +3.12%  "
    12 = "  +1.11%  "
    13 = "  -0.25%  "
    14 = "  -0.92%  "
    15 = "  -2.27%  "
    16 = "  +0.45%  "
    17 = "  +0.38%  "
    18 = "  -0.81%  "
    19 = "  +0.65%  "
    20 = "  +0.96%  "
    21 = "  -0.33%  "
    22 = "  -0.05%  "
    23 = "  +0.64%  "
    24 = "  +2.99%  "
    25 = "  +5.38%  "
    26 = "  +2.85%  "
    27 = "  -1.04%  "
    28 = "  -7.86%  "
    29 = "  +5.00%  "
    30 = "  +0.07%  "
    31 = "  -0.24%  "
    32 = "  +0.67%  "
    33 = "  +1.42%  "
    34 = "  +1.24%  "
    35 = "  -0.41%  "
    37 = "  +0.50%  "
    38 = "  -0.53%  "
    39 = "  +1.08%  "
    41 = "  +0.50%  "
    42 = "  +0.06%  "
    43 = "  +0.49%  "
    44 = "  +0.02%  "
    45 = "  -0.41%  "
    46 = "  +0.70%  "
    47 = "  +0.51%  "
    48 = "  +0.13%  "
    49 = "  +0.29%  "
    50 = "  +0.37%  "
    51 = "  +0.51%  "
}

# Price strings that Excel's normal numeric auto-detection would leave
# alone (they contain two "." separators, e.g. "1.687.86", so they can
# never be parsed as a single number and stay text on their own).
$naturallyTextRows = @(2, 3, 12, 17, 39)

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    if ($naturallyTextRows -notcontains $row) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $priceUpdates[$row]
}

foreach ($row in $volumeUpdates.Keys) {
    $ws.Cells.Item($row, 5).Value = $volumeUpdates[$row]
}
